# Auto-generated Excel COM-interop edit script
# Applies the numeric updates described in the commit diff to the
# "Sheets" workbook (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR tabs).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H74").Value = 13100
$ws.Range("I74").Value = 5166.6665
$ws.Range("K74").Value = 5166.6665
$ws.Range("M74").Value = -4230.6665
$ws.Range("H77").Value = 13100
$ws.Range("I77").Value = 5166.6665
$ws.Range("K77").Value = 25833.3325
$ws.Range("M77").Value = -21153.3325
$ws.Range("H98").Value = 1361.875
$ws.Range("I98").Value = 1361.875
$ws.Range("K98").Value = 1361.875
$ws.Range("M98").Value = 136.125
$ws.Range("H107").Value = 132.13333
$ws.Range("I107").Value = 75.76922999999999
$ws.Range("J107").Value = 498.5
$ws.Range("K107").Value = 75.76922999999999
$ws.Range("L107").Value = 498.5
$ws.Range("M107").Value = 1844.23077
$ws.Range("N107").Value = -4338.5
$ws.Range("H113").Value = 2000
$ws.Range("J113").Value = 2000
$ws.Range("L113").Value = 2000
$ws.Range("N113").Value = -8508
$ws.Range("H115").Value = 783.53845
$ws.Range("I115").Value = 765.5
$ws.Range("J115").Value = 1000
$ws.Range("K115").Value = 2296.5
$ws.Range("L115").Value = 3000
$ws.Range("M115").Value = -729.5
$ws.Range("N115").Value = -6134
$ws.Range("H122").Value = 1361.875
$ws.Range("I122").Value = 1361.875
$ws.Range("K122").Value = 4085.625
$ws.Range("M122").Value = -1635.625
$ws.Range("H127").Value = 1934.9
$ws.Range("J127").Value = 1108.5
$ws.Range("L127").Value = 3325.5
$ws.Range("N127").Value = -13245.5
$ws.Range("H131").Value = 1411.25
$ws.Range("I131").Value = 748.25
$ws.Range("J131").Value = 2074.25
$ws.Range("K131").Value = 2244.75
$ws.Range("L131").Value = 6222.75
$ws.Range("M131").Value = 2795.25
$ws.Range("N131").Value = -16302.75
$ws.Range("H138").Value = 1779.2106
$ws.Range("I138").Value = 1507.6923
$ws.Range("J138").Value = 2367.5
$ws.Range("K138").Value = 4523.0769
$ws.Range("L138").Value = 7102.5
$ws.Range("M138").Value = 616.9231
$ws.Range("N138").Value = -17382.5

$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H61").Value = 7743.75
$ws.Range("I61").Value = 7743.75
$ws.Range("K61").Value = 7743.75
$ws.Range("M61").Value = -7531.75
$ws.Range("H132").Value = 3179.6
$ws.Range("I132").Value = 1828.1428
$ws.Range("K132").Value = 5484.428400000001
$ws.Range("M132").Value = -2954.428400000001
$ws.Range("H136").Value = 7743.75
$ws.Range("I136").Value = 7743.75
$ws.Range("K136").Value = 23231.25
$ws.Range("M136").Value = -20681.25

$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H86").Value = 6036.3335
$ws.Range("I86").Value = 488.6
$ws.Range("K86").Value = 488.6
$ws.Range("M86").Value = 634.4
$ws.Range("H89").Value = 6036.3335
$ws.Range("I89").Value = 488.6
$ws.Range("K89").Value = 2443
$ws.Range("M89").Value = 3173
$ws.Range("H107").Value = 4181.048
$ws.Range("I107").Value = 1292.4615
$ws.Range("K107").Value = 1292.4615
$ws.Range("M107").Value = 627.5385000000001

$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H31").Value = 6241.615
$ws.Range("I31").Value = 1455.091
$ws.Range("K31").Value = 1455.091
$ws.Range("M31").Value = -1160.091
$ws.Range("H34").Value = 6241.615
$ws.Range("I34").Value = 1455.091
$ws.Range("K34").Value = 1455.091
$ws.Range("M34").Value = -1253.091
$ws.Range("H94").Value = 3249.5715
$ws.Range("J94").Value = 5172.7144
$ws.Range("L94").Value = 5172.7144
$ws.Range("N94").Value = -6074.7144
$ws.Range("H122").Value = 1488.9524
$ws.Range("J122").Value = 2236
$ws.Range("L122").Value = 6708
$ws.Range("N122").Value = -11608

$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H63").Value = 1875.5
$ws.Range("I63").Value = 1875.5
$ws.Range("K63").Value = 5626.5
$ws.Range("M63").Value = -4877.5
$ws.Range("H66").Value = 1875.5
$ws.Range("I66").Value = 1875.5
$ws.Range("K66").Value = 16879.5
$ws.Range("M66").Value = -13135.5
$ws.Range("H86").Value = 207.8
$ws.Range("I86").Value = 195.75
$ws.Range("K86").Value = 587.25
$ws.Range("M86").Value = 598.75
$ws.Range("H89").Value = 207.8
$ws.Range("I89").Value = 195.75
$ws.Range("K89").Value = 1761.75
$ws.Range("M89").Value = 4166.25
$ws.Range("H107").Value = 868.9231
$ws.Range("I107").Value = 701.5
$ws.Range("J107").Value = 899.36365
$ws.Range("K107").Value = 2104.5
$ws.Range("L107").Value = 2698.09095
$ws.Range("M107").Value = -184.5
$ws.Range("N107").Value = -6538.09095
$ws.Range("H113").Value = 2481.3333
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 2481.3333
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 7443.999899999999
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -11783.9999
$ws.Range("H129").Value = 1296.7142
$ws.Range("I129").Value = 855.4
$ws.Range("J129").Value = 2400
$ws.Range("K129").Value = 2566.2
$ws.Range("L129").Value = 7200
$ws.Range("M129").Value = 2433.8
$ws.Range("N129").Value = -17200

$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H102").Value = 1248.9706
$ws.Range("I102").Value = 886.4194
$ws.Range("K102").Value = 886.4194
$ws.Range("M102").Value = 735.5806
$ws.Range("H113").Value = 8138.9
$ws.Range("J113").Value = 8736.25
$ws.Range("L113").Value = 8736.25
$ws.Range("N113").Value = -13076.25
$ws.Range("H122").Value = 141448.14
$ws.Range("I122").Value = 201929.95
$ws.Range("J122").Value = 3989.4546
$ws.Range("K122").Value = 605789.8500000001
$ws.Range("L122").Value = 11968.3638
$ws.Range("M122").Value = -603339.8500000001
$ws.Range("N122").Value = -16868.3638
$ws.Range("H126").Value = 3877.9333
$ws.Range("I126").Value = 3472.5833
$ws.Range("J126").Value = 5499.3335
$ws.Range("K126").Value = 10417.7499
$ws.Range("L126").Value = 16498.0005
$ws.Range("M126").Value = -7947.749899999999
$ws.Range("N126").Value = -21438.0005

$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H122").Value = 4123.375
$ws.Range("I122").Value = 3831.3333
$ws.Range("J122").Value = 4999.5
$ws.Range("K122").Value = 11493.9999
$ws.Range("L122").Value = 14998.5
$ws.Range("M122").Value = -9043.999899999999
$ws.Range("N122").Value = -19898.5
$ws.Range("H132").Value = 4179
$ws.Range("I132").Value = 3406.3333
$ws.Range("K132").Value = 10218.9999
$ws.Range("M132").Value = -7688.999899999999

$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H122").Value = 1465.7778
$ws.Range("I122").Value = 1493
$ws.Range("K122").Value = 4479
$ws.Range("M122").Value = -2029
$ws.Range("H126").Value = 2954.2
$ws.Range("I126").Value = 1446
$ws.Range("J126").Value = 5216.5
$ws.Range("K126").Value = 4338
$ws.Range("L126").Value = 15649.5
$ws.Range("M126").Value = -1868
$ws.Range("N126").Value = -20589.5
$ws.Range("H136").Value = 2926.9062
$ws.Range("J136").Value = 6965.2856
$ws.Range("L136").Value = 20895.8568
$ws.Range("N136").Value = -25995.8568
